# The presentation ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (the "live" theme)
#   ppt/theme/theme2.xml  -> bound to the notes master
#
# The target edit swaps their contents: theme1.xml becomes the stock
# "Office Theme" palette (previously living in theme2.xml) and
# theme2.xml becomes the "Integral" palette (previously living in
# theme1.xml). Both themes already share an identical font scheme
# (Arial/Arial) and format scheme, so only the 10 color-scheme entries
# that actually differ between the two themes need to change.
#
# PowerPoint's ThemeColorScheme indexes the 12 theme colors in the
# fixed OOXML clrScheme order: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink. RGB() isn't available in this host, so the BGR long
# values (0x00BBGGRR) are written out explicitly.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeElements.ThemeColorScheme

# Target palette: "Office Theme" (swapped into theme1.xml)
$colorScheme.Item(1).RGB  = 0x000000   # dk1      000000
$colorScheme.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 0x6A5444   # dk2      44546A
$colorScheme.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$colorScheme.Item(9).RGB  = 0xC47244   # accent5  4472C4
$colorScheme.Item(10).RGB = 0x47AD70   # accent6  70AD47
$colorScheme.Item(11).RGB = 0xC16305   # hlink    0563C1
$colorScheme.Item(12).RGB = 0x724F95   # folHlink 954F72
